$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = 0
$ws.Cells.Item(2, 3).Value = 0
$ws.Cells.Item(3, 2).Value = 0
$ws.Cells.Item(3, 3).Value = 0
$ws.Cells.Item(4, 2).Value = 0
$ws.Cells.Item(4, 3).Value = 0
$ws.Cells.Item(5, 2).Value = 0
$ws.Cells.Item(5, 3).Value = 0
$ws.Cells.Item(6, 2).Value = 0
$ws.Cells.Item(6, 3).Value = 0
$ws.Cells.Item(7, 2).Value = 0
$ws.Cells.Item(7, 3).Value = 0
$ws.Cells.Item(8, 2).Value = 0
$ws.Cells.Item(8, 3).Value = 0
$ws.Cells.Item(9, 2).Value = 0
$ws.Cells.Item(9, 3).Value = 0
$ws.Cells.Item(10, 2).Value = 0
$ws.Cells.Item(10, 3).Value = 0
$ws.Cells.Item(11, 2).Value = 0
$ws.Cells.Item(11, 3).Value = 0
$ws.Cells.Item(12, 2).Value = 0
$ws.Cells.Item(12, 3).Value = 0
$ws.Cells.Item(13, 2).Value = 0
$ws.Cells.Item(13, 3).Value = 0
$ws.Cells.Item(14, 2).Value = 0.045353539288043976
$ws.Cells.Item(14, 3).Value = 0.070209339261054993
$ws.Cells.Item(15, 2).Value = 0.54213064908981323
$ws.Cells.Item(15, 3).Value = 0.24002590775489807
$ws.Cells.Item(16, 2).Value = 1.8167499303817749
$ws.Cells.Item(16, 3).Value = 0.28282567858695984
$ws.Cells.Item(17, 2).Value = 3.5960521697998047
$ws.Cells.Item(17, 3).Value = 0.62914460897445679
$ws.Cells.Item(18, 2).Value = 5.3597936630249023
$ws.Cells.Item(18, 3).Value = 1.2266266345977783
$ws.Cells.Item(19, 2).Value = 6.603663444519043
$ws.Cells.Item(19, 3).Value = 1.7457572221755981
$ws.Cells.Item(20, 2).Value = 7.0056796073913574
$ws.Cells.Item(20, 3).Value = 1.9818108081817627
$ws.Cells.Item(21, 2).Value = 6.4986319541931152
$ws.Cells.Item(21, 3).Value = 1.8624374866485596
$ws.Cells.Item(22, 2).Value = 5.2605199813842773
$ws.Cells.Item(22, 3).Value = 1.4526060819625854
$ws.Cells.Item(23, 2).Value = 3.6424412727355957
$ws.Cells.Item(23, 3).Value = 0.93047523498535156
$ws.Cells.Item(24, 2).Value = 2.0579819679260254
$ws.Cells.Item(24, 3).Value = 0.53905999660491943
$ws.Cells.Item(25, 2).Value = 0.85927265882492065
$ws.Cells.Item(25, 3).Value = 0.39832442998886108
$ws.Cells.Item(26, 2).Value = 0.22252239286899567
$ws.Cells.Item(26, 3).Value = 0.27227833867073059
$ws.Cells.Item(27, 2).Value = 0.06060342863202095
$ws.Cells.Item(27, 3).Value = 0.083006061613559723
$ws.Cells.Item(28, 2).Value = 0.027220681309700012
$ws.Cells.Item(28, 3).Value = 0.060867298394441605
$ws.Cells.Item(29, 2).Value = 0.040450219064950943
$ws.Cells.Item(29, 3).Value = 0.090449444949626923
$ws.Cells.Item(30, 2).Value = 0.056257877498865128
$ws.Cells.Item(30, 3).Value = 0.12579643726348877
$ws.Cells.Item(31, 2).Value = 0.1329600065946579
$ws.Cells.Item(31, 3).Value = 0.12701381742954254
$ws.Cells.Item(32, 2).Value = 0.31399592757225037
$ws.Cells.Item(32, 3).Value = 0.17910207808017731
$ws.Cells.Item(33, 2).Value = 0.5460735559463501
$ws.Cells.Item(33, 3).Value = 0.32781586050987244
$ws.Cells.Item(34, 2).Value = 0.74393594264984131
$ws.Cells.Item(34, 3).Value = 0.47304326295852661
$ws.Cells.Item(35, 2).Value = 0.83866149187088013
$ws.Cells.Item(35, 3).Value = 0.5527070164680481
$ws.Cells.Item(36, 2).Value = 0.80513209104537964
$ws.Cells.Item(36, 3).Value = 0.55073446035385132
$ws.Cells.Item(37, 2).Value = 0.66804289817810059
$ws.Cells.Item(37, 3).Value = 0.4873606264591217
$ws.Cells.Item(38, 2).Value = 0.48890575766563416
$ws.Cells.Item(38, 3).Value = 0.40287584066390991
$ws.Cells.Item(39, 2).Value = 0.33851951360702515
$ws.Cells.Item(39, 3).Value = 0.34200698137283325
$ws.Cells.Item(40, 2).Value = 0.26021048426628113
$ws.Cells.Item(40, 3).Value = 0.34244221448898315
$ws.Cells.Item(41, 2).Value = 0.22890013456344604
$ws.Cells.Item(41, 3).Value = 0.40367940068244934
$ws.Cells.Item(42, 2).Value = 0.21565344929695129
$ws.Cells.Item(42, 3).Value = 0.45521152019500732
$ws.Cells.Item(43, 2).Value = 0.20387516915798187
$ws.Cells.Item(43, 3).Value = 0.45587876439094543
$ws.Cells.Item(44, 2).Value = 0.1809966117143631
$ws.Cells.Item(44, 3).Value = 0.40472075343132019
$ws.Cells.Item(45, 2).Value = 0.14243166148662567
$ws.Cells.Item(45, 3).Value = 0.31848686933517456
$ws.Cells.Item(46, 2).Value = 0.096858836710453033
$ws.Cells.Item(46, 3).Value = 0.21658295392990112
$ws.Cells.Item(47, 2).Value = 0.054015733301639557
$ws.Cells.Item(47, 3).Value = 0.12078285217285156
$ws.Cells.Item(48, 2).Value = 0.021957110613584518
$ws.Cells.Item(48, 3).Value = 0.049097590148448944
$ws.Cells.Item(49, 2).Value = 0.01248572114855051
$ws.Cells.Item(49, 3).Value = 0.018074385821819305
$ws.Cells.Item(50, 2).Value = 0.027247501537203789
$ws.Cells.Item(50, 3).Value = 0.060927268117666245
$ws.Cells.Item(51, 2).Value = 0.050809253007173538
$ws.Cells.Item(51, 3).Value = 0.11361294984817505
$ws.Cells.Item(52, 2).Value = 0.068109974265098572
$ws.Cells.Item(52, 3).Value = 0.15229853987693787
$ws.Cells.Item(53, 2).Value = 0.073511689901351929
$ws.Cells.Item(53, 3).Value = 0.15823647379875183
$ws.Cells.Item(54, 2).Value = 0.077044554054737091
$ws.Cells.Item(54, 3).Value = 0.12665265798568726
$ws.Cells.Item(55, 2).Value = 0.10914053022861481
$ws.Cells.Item(55, 3).Value = 0.11959797888994217
$ws.Cells.Item(56, 2).Value = 0.22455909848213196
$ws.Cells.Item(56, 3).Value = 0.24478836357593536
$ws.Cells.Item(57, 2).Value = 0.49640658497810364
$ws.Cells.Item(57, 3).Value = 0.43962976336479187
$ws.Cells.Item(58, 2).Value = 0.90357577800750732
$ws.Cells.Item(58, 3).Value = 0.68130695819854736
$ws.Cells.Item(59, 2).Value = 1.3401745557785034
$ws.Cells.Item(59, 3).Value = 0.98103153705596924
$ws.Cells.Item(60, 2).Value = 1.6949303150177002
$ws.Cells.Item(60, 3).Value = 1.3277708292007446
$ws.Cells.Item(61, 2).Value = 1.8814398050308228
$ws.Cells.Item(61, 3).Value = 1.7079257965087891
$ws.Cells.Item(62, 2).Value = 1.86708402633667
$ws.Cells.Item(62, 3).Value = 2.0817947387695313
$ws.Cells.Item(63, 2).Value = 1.6824629306793213
$ws.Cells.Item(63, 3).Value = 2.3649694919586182
$ws.Cells.Item(64, 2).Value = 1.4062883853912354
$ws.Cells.Item(64, 3).Value = 2.4514586925506592
$ws.Cells.Item(65, 2).Value = 1.1213099956512451
$ws.Cells.Item(65, 3).Value = 2.2722933292388916
$ws.Cells.Item(66, 2).Value = 0.837727427482605
$ws.Cells.Item(66, 3).Value = 1.8533576726913452
$ws.Cells.Item(67, 2).Value = 0.567320704460144
$ws.Cells.Item(67, 3).Value = 1.2685675621032715
$ws.Cells.Item(68, 2).Value = 0.32230293750762939
$ws.Cells.Item(68, 3).Value = 0.66315603256225586
$ws.Cells.Item(69, 2).Value = 0.17106030881404877
$ws.Cells.Item(69, 3).Value = 0.237416610121727
$ws.Cells.Item(70, 2).Value = 0.14556881785392761
$ws.Cells.Item(70, 3).Value = 0.32550176978111267
$ws.Cells.Item(71, 2).Value = 0.20812630653381348
$ws.Cells.Item(71, 3).Value = 0.46538460254669189
$ws.Cells.Item(72, 2).Value = 0.23801761865615845
$ws.Cells.Item(72, 3).Value = 0.53222358226776123
